# hed3_tags_single_sheet.xlsx - "Update HedString further, use in more places.
# Update validation to only work on long tags.  Misc other cleanup and fixes"
#
# Concrete spreadsheet-level changes made by this commit:
#   1. Fix a typo in D2: "Attribute/Sensory/Bisual" -> "Attribute/Sensory/Visual".
#   2. Replace the long HED tag string in D3 with the shorter "Awake".
#   3. Remove the two example rows (4 and 5) entirely - ResponseOnset /
#      ResponseOffset - but keep their row heights around (row 4 stays
#      47pt tall, row 5 settles back down to ~13.8pt).
#   4. Row 3 no longer needs the taller wrapped-text height, so it shrinks
#      from 28.5pt to 14.9pt.
#   5. The active selection moves from D3 to A4, and the view scrolls back
#      to show column A (topLeftCell A1 instead of D1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "Bisual" typo.
$ws.Range("D2").Value = "Attribute/Sensory/Visual"

# 2. Shorten the D3 HED tag string down to "Awake".
$ws.Range("D3").Value = "Awake"

# 3. Delete the contents of the two trailing example rows (ResponseOnset /
#    ResponseOffset), leaving the now-blank rows in place.
$ws.Range("A4:D5").Clear()

# 4. Adjust row heights: row 3 shrinks now that D3 is a short string, row 4
#    keeps its existing (taller) height, row 5 settles to its natural height.
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(4).RowHeight = 47
$ws.Rows.Item(5).RowHeight = 13.8

# 5. Move the selection/view to A4 (first of the now-empty rows), scrolled
#    back to column A.
$ws.Range("A4").Select() | Out-Null
